# Aggiornamento 15, 16, 17 marzo
# Adds three new rows (227, 228, 229) of data to the end of the sheet,
# continuing the existing table of "nuovi pos." / rolling 7-day sum /
# rolling 7-day sum per 100k inhabitants.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$newRows = @(
    @(44301, 0, 2, 93.41429238673517),
    @(44302, 0, 2, 93.41429238673517),
    @(44303, 1, 1, 46.70714619336758)
)

$lastExistingRow = 226
$startRow = 227

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    # Copy the formatting (number format, font, border, alignment) from the
    # last existing row's date cell so the new date cell matches the style
    # of the rest of the column.
    $ws.Cells.Item($lastExistingRow, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}

$excel.CutCopyMode = $false
